# Generate Report for Handback
#
# Refresh the localization handback-status report: the handoff/handback
# timestamps for the "3537563b-...md" source file (row 2 on each
# language sheet) were regenerated, which also updates the roll-up
# "Latest HO Xliff Generate Date" cell for that row on the Overview
# sheet. The f5b1a7b0-...md row (row 3) is untouched.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 2 = 3537563b-...md
#   H = Correspond Handoff Datetime, K = Correspond Handback DateTime
$zhcn.Range("H2").Value = "2016-08-13 23:03:49"
$zhcn.Range("K2").Value = "2016-08-13 23:04:17"

# de-de sheet: row 2 = 3537563b-...md
$dede.Range("H2").Value = "2016-08-13 23:03:56"
$dede.Range("K2").Value = "2016-08-13 23:04:27"

# Overview sheet: row 2 = 3537563b-...md, column G = Latest HO Xliff Generate Date
$overview.Range("G2").Value = "2016-08-13 23:03:56"
